# Daily attendance processing - reorder "Recorded By" (column G) entries so
# that the "System" contributor is listed first, followed by the remaining
# contributors in alphabetical order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 157

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $raw = $cell.Text

    if ($raw -eq $null -or $raw -eq "") {
        continue
    }

    $parts = $raw -split ","
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $systemParts = @()
    $otherParts = @()
    foreach ($p in $trimmed) {
        if ($p.Equals("System")) {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    $otherSorted = $otherParts | Sort-Object
    $ordered = $systemParts + $otherSorted
    $newValue = [string]::Join(", ", $ordered)

    if ($newValue -ne $raw) {
        $cell.Value = $newValue
    }
}
